# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker data table (B16:J19) is reorganized from being grouped by
# worker (each worker's two overdue periods on consecutive rows) to being
# grouped by period (each period's two workers on consecutive rows).
#
# Resulting layout:
#   Row16: CC 73203421  ORLANDO ANTONIO LORDUY FLOREZ  period 1705  120000 / 3000000
#   Row17: CC 1047480732 FREDDY JOHANES VARGAS RAMIREZ period 1705  36000  / 900000
#   Row18: CC 73203421  ORLANDO ANTONIO LORDUY FLOREZ  period 1706  120000 / 3000000
#   Row19: CC 1047480732 FREDDY JOHANES VARGAS RAMIREZ period 1706  36000  / 900000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing contents of the data rows (keeps existing cell formatting/borders intact)
$ws.Range("B16:J19").ClearContents()

# Row 16: Orlando / period 1705
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73203421"
$ws.Range("D16").Value = "ORLANDO ANTONIO LORDUY FLOREZ"
$ws.Range("E16").Value = "1705"
$ws.Range("F16").Value = 120000
$ws.Range("G16").Value = 3000000

# Row 17: Freddy / period 1705
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047480732"
$ws.Range("D17").Value = "FREDDY JOHANES VARGAS RAMIREZ"
$ws.Range("E17").Value = "1705"
$ws.Range("F17").Value = 36000
$ws.Range("G17").Value = 900000

# Row 18: Orlando / period 1706
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73203421"
$ws.Range("D18").Value = "ORLANDO ANTONIO LORDUY FLOREZ"
$ws.Range("E18").Value = "1706"
$ws.Range("F18").Value = 120000
$ws.Range("G18").Value = 3000000

# Row 19: Freddy / period 1706
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047480732"
$ws.Range("D19").Value = "FREDDY JOHANES VARGAS RAMIREZ"
$ws.Range("E19").Value = "1706"
$ws.Range("F19").Value = 36000
$ws.Range("G19").Value = 900000
